# Update attendee/view count figures (column F) on the "展览" and "全部类型"
# sheets to reflect newly generated data (output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of old value -> new value for column F on each affected row.
$cellUpdates = @{
    "F2"  = 293
    "F3"  = 1178
    "F4"  = 16713
    "F6"  = 1634
    "F9"  = 369
    "F12" = 11605
    "F14" = 1279
    "F15" = 4592
    "F16" = 423
    "F19" = 884
    "F21" = 151
}

# "全部类型" uses a different row numbering (it includes extra rows for
# other event types interspersed), so cells are addressed separately.
$cellUpdatesAll = @{
    "F2"  = 293
    "F4"  = 1178
    "F5"  = 16713
    "F7"  = 1634
    "F10" = 369
    "F15" = 11605
    "F17" = 1279
    "F18" = 4592
    "F19" = 423
    "F22" = 884
    "F24" = 151
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($addr in $cellUpdates.Keys) {
    $ws1.Range($addr).Value = $cellUpdates[$addr]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($addr in $cellUpdatesAll.Keys) {
    $ws4.Range($addr).Value = $cellUpdatesAll[$addr]
}
